$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Extend header row (row 1) with two new columns P and Q ---
# Copy formatting from the existing header cell O1 (style s="1") onto P1 and Q1
# so the new cells render with the same bold/border/centered style as the rest
# of row 1, then set their values.
$ws.Range("O1").Copy() | Out-Null
$ws.Range("P1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("Q1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# --- Update data rows 2-25 ---
# Columns I, K, M, O swap their 1/2 values, and two new columns P, Q (value 2)
# are appended to every row.
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value = 2   # I: 1 -> 2
    $ws.Cells.Item($r, 11).Value = 1  # K: 2 -> 1
    $ws.Cells.Item($r, 13).Value = 2  # M: 1 -> 2
    $ws.Cells.Item($r, 15).Value = 1  # O: 2 -> 1
    $ws.Cells.Item($r, 16).Value = 2  # P: new column
    $ws.Cells.Item($r, 17).Value = 2  # Q: new column
}
